$wb = $excel.ActiveWorkbook

# Sheets that contain the mirrored "展览" data table: "展览" (1) and "全部类型" (4)
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- 1) Bump "想去人数" (column F) counters on rows 2-24 that only changed their count ---
    $ws.Cells.Item(3, 6).Value = 3069
    $ws.Cells.Item(4, 6).Value = 224
    $ws.Cells.Item(6, 6).Value = 196
    $ws.Cells.Item(7, 6).Value = 1655
    $ws.Cells.Item(8, 6).Value = 1618
    $ws.Cells.Item(9, 6).Value = 54
    $ws.Cells.Item(10, 6).Value = 358
    $ws.Cells.Item(14, 6).Value = 28
    $ws.Cells.Item(15, 6).Value = 226
    $ws.Cells.Item(20, 6).Value = 42
    $ws.Cells.Item(21, 6).Value = 13
    $ws.Cells.Item(22, 6).Value = 363
    $ws.Cells.Item(23, 6).Value = 172

    # --- 2) Insert a brand-new row 25 for the new event, shifting old rows 25-36 to 26-37 ---
    $ws.Rows.Item(25).Insert()

    # Restore the bordered/bold/centered look ("style 1") on the new index cell, same as
    # every other cell in column A of the data table.
    $newIndexCell = $ws.Cells.Item(25, 1)
    $newIndexCell.Font.Bold = $true
    $newIndexCell.HorizontalAlignment = -4108
    $newIndexCell.VerticalAlignment = -4160
    $newIndexCell.Borders.LineStyle = 1

    # --- 3) Populate the new row 25 with the new event's data ---
    # A leading apostrophe forces the YYYY-MM-DD-looking date string to be
    # stored as plain text instead of Excel auto-converting it to a date.
    $ws.Cells.Item(25, 1).Value = 24
    $ws.Cells.Item(25, 2).Value = "'2024-07-21"
    $ws.Cells.Item(25, 3).Value = "乐平·CY境界次元动漫夏时庆"
    $ws.Cells.Item(25, 4).Value = "翥山西路182号 佳佳基大酒店"
    $ws.Cells.Item(25, 5).Value = "2024.07.21 10:00-07.21 17:00"
    $ws.Cells.Item(25, 6).Value = 0
    $ws.Cells.Item(25, 7).Value = 30
    $ws.Cells.Item(25, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86768"
    $ws.Cells.Item(25, 9).Value = "//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png"

    # --- 4) Renumber the sequential index in column A for every row pushed down by the insert ---
    for ($r = 26; $r -le 37; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # --- 5) Apply the extra "想去人数" bumps that landed on the shifted rows ---
    $ws.Cells.Item(27, 6).Value = 2050
    $ws.Cells.Item(31, 6).Value = 189
    $ws.Cells.Item(36, 6).Value = 498
    $ws.Cells.Item(37, 6).Value = 7
}
